$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Emprestimos")

# New header labels (variation columns)
$ws.Range("G4").Value = "3T25/2T25"
$ws.Range("H4").Value = "3T25/3T24"

# New variation values
$ws.Range("G5").Value = 0.223
$ws.Range("H5").Value = -0.334

$ws.Range("G6").Value = -0.03
$ws.Range("H6").Value = -0.11

$ws.Range("G7").Value = -0.136
$ws.Range("H7").Value = -0.248

$ws.Range("G8").Value = 0.063
$ws.Range("H8").Value = 0.253

$ws.Range("G9").Value = 0.008
$ws.Range("H9").Value = 0.037

$ws.Range("G10").Value = -0.154
$ws.Range("H10").Value = -0.266

# Percentage number format for the new variation columns
$ws.Range("G5:H10").Style = "Percent"
$ws.Range("G5:H10").NumberFormat = "0.0%"

# Column widths for the new columns
$ws.Columns.Item("G").ColumnWidth = 9.7109375
$ws.Columns.Item("H").ColumnWidth = 9.7109375

# Selection / active cell on this sheet
$ws.Range("N7").Select()

# Make "Emprestimos" the active (selected) sheet/tab
$ws.Activate()

# The previously active sheet ("Seguros e Cartoes") keeps its own last
# selection, unaffected.
